{"js": "// Applies the \"Major changes relating to itinerary master\" edit:\n//   - \"Itinerary_no-Itinerary_sub_no-Version-ddmmyyyy\" -> \"Number-SubNumber-DayNumber-Version\"\n//   - \"Itinerary_no is unique ...\" -> \"Number is unique ...\"\n//   - \"Itinerary_sub_no is by default 1 ...\" -> \"SubNumber is by default 1 ...\"\n//   - \"ddmmyyyy will be date of creation\" -> \"DayNumber will have values 1,2,3\u2026 representing Day 1, Day 2, Day 3 of the itinerary\"\n//   - \"-01222015\" -> \"-5\" (inside the \"(Example: 1234-1-2-01222015)\" sample)\n//\n// (The remaining hunks in the source diff only wrap existing, unchanged\n// words \u2014 \"mico\", \"reporting\", \"Spring\", \"pdf\" \u2014 with <w:proofErr/> spell/\n// grammar-check markers that Word's proofing engine stamps in as a side\n// effect of the edit session; they carry no visible text change, so there\n// is nothing further to apply for those spans.)\n\nconst body = context.document.body;\n\nasync function replaceOnce(findText, newText) {\n  const results = body.search(findText, { matchCase: true, matchWildcards: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + findText);\n  }\n\n  // Every target string below is unique in the document; replace each\n  // occurrence found (in practice exactly one) so the script is resilient\n  // even if a string were duplicated elsewhere.\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1) Itinerary Format line.\nawait replaceOnce(\n  \"Itinerary_no-Itinerary_sub_no-Version-ddmmyyyy\",\n  \"Number-SubNumber-DayNumber-Version\"\n);\n\n// 2) Itinerary_no -> Number definition.\nawait replaceOnce(\n  \"Itinerary_no is unique and incremental. It will remain unique across all years and never repeated. Any itinerary may be reused but will be assigned a unique number once it is created.\",\n  \"Number is unique and incremental. It will remain unique across all years and never repeated. Any itinerary may be reused but will be assigned a unique number once it is created.\"\n);\n\n// 3) Itinerary_sub_no -> SubNumber definition.\nawait replaceOnce(\n  \"Itinerary_sub_no is by default 1 meaning every Itinerary will have a minimum of one sub itinerary. In case, if multiple itineraries are to be maintained for larger groups, then sub itineraries can handle this special requirement while being part of main itinerary. For example some inbound members of a group may want to visit Chandigarh while visiting North India and others some other city, say Shimla. This will help manage this requirement.\",\n  \"SubNumber is by default 1 meaning every Itinerary will have a minimum of one sub itinerary. In case, if multiple itineraries are to be maintained for larger groups, then sub itineraries can handle this special requirement while being part of main itinerary. For example some inbound members of a group may want to visit Chandigarh while visiting North India and others some other city, say Shimla. This will help manage this requirement.\"\n);\n\n// 4) ddmmyyyy -> DayNumber definition (and rewritten explanation).\nawait replaceOnce(\n  \"ddmmyyyy will be date of creation\",\n  \"DayNumber will have values 1,2,3\\u2026 representing Day 1, Day 2, Day 3 of the itinerary\"\n);\n\n// 5) Example date suffix -> day number.\nawait replaceOnce(\n  \"-01222015\",\n  \"-5\"\n);\n", "ps1": "# Applies the \"Major changes relating to itinerary master\" edit:\n#   - \"Itinerary_no-Itinerary_sub_no-Version-ddmmyyyy\" -> \"Number-SubNumber-DayNumber-Version\"\n#   - \"Itinerary_no is unique ...\" -> \"Number is unique ...\"\n#   - \"Itinerary_sub_no is by default 1 ...\" -> \"SubNumber is by default 1 ...\"\n#   - \"ddmmyyyy will be date of creation\" -> \"DayNumber will have values 1,2,3\u2026 representing Day 1, Day 2, Day 3 of the itinerary\"\n#   - \"-01222015\" -> \"-5\" (inside the \"(Example: 1234-1-2-01222015)\" sample)\n#\n# (The remaining hunks in the source diff only wrap existing, unchanged\n# words -- \"mico\", \"reporting\", \"Spring\", \"pdf\" -- with <w:proofErr/> spell/\n# grammar-check markers that Word's proofing engine stamps in as a side\n# effect of the edit session; they carry no visible text change, so there\n# is nothing further to apply for those spans.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once($findText, $replaceText) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.MatchCase = $true\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $found = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    if (-not $found) {\n        throw \"Find/Replace did not find text: $findText\"\n    }\n}\n\n# 1) Itinerary Format line.\nReplace-Once \"Itinerary_no-Itinerary_sub_no-Version-ddmmyyyy\" \"Number-SubNumber-DayNumber-Version\"\n\n# 2) Itinerary_no -> Number definition.\nReplace-Once \"Itinerary_no is unique and incremental. It will remain unique across all years and never repeated. Any itinerary may be reused but will be assigned a unique number once it is created.\" \"Number is unique and incremental. It will remain unique across all years and never repeated. Any itinerary may be reused but will be assigned a unique number once it is created.\"\n\n# 3) Itinerary_sub_no -> SubNumber definition.\nReplace-Once \"Itinerary_sub_no is by default 1 meaning every Itinerary will have a minimum of one sub itinerary. In case, if multiple itineraries are to be maintained for larger groups, then sub itineraries can handle this special requirement while being part of main itinerary. For example some inbound members of a group may want to visit Chandigarh while visiting North India and others some other city, say Shimla. This will help manage this requirement.\" \"SubNumber is by default 1 meaning every Itinerary will have a minimum of one sub itinerary. In case, if multiple itineraries are to be maintained for larger groups, then sub itineraries can handle this special requirement while being part of main itinerary. For example some inbound members of a group may want to visit Chandigarh while visiting North India and others some other city, say Shimla. This will help manage this requirement.\"\n\n# 4) ddmmyyyy -> DayNumber definition (and rewritten explanation).\nReplace-Once \"ddmmyyyy will be date of creation\" \"DayNumber will have values 1,2,3\u2026 representing Day 1, Day 2, Day 3 of the itinerary\"\n\n# 5) Example date suffix -> day number.\nReplace-Once \"-01222015\" \"-5\"\n"}
